# The deck currently ships the "Integral" theme as the live/applied theme
# (colour scheme: dk2=455F51, accent1=99CB38, accent2=63A537, accent3=E6D024,
# accent4=CC9700, accent5=4EB3CF, accent6=378DA6, hlink=6B9F25,
# folHlink=B26B02) while the stock "Office Theme" palette
# (dk2=44546A, accent1=5B9BD5, accent2=ED7D31, accent3=A5A5A5, accent4=FFC000,
# accent5=4472C4, accent6=70AD47, hlink=0563C1, folHlink=954F72) sits unused
# in the companion theme part. The edit swaps which palette is live, i.e.
# re-colours the presentation's active theme back to the default Office
# colours.
#
# Theme colours are addressed through Slide.ThemeColorScheme, whose 12 slots
# line up with <a:clrScheme> in document order: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink. RGB values use the usual COM 0xBBGGRR packing.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function Set-ThemeColor($scheme, $index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $scheme.Item($index).RGB = $r -bor ($g -shl 8) -bor ($b -shl 16)
}

Set-ThemeColor $tcs 1  "000000"  # dk1
Set-ThemeColor $tcs 2  "FFFFFF"  # lt1
Set-ThemeColor $tcs 3  "44546A"  # dk2
Set-ThemeColor $tcs 4  "E7E6E6"  # lt2
Set-ThemeColor $tcs 5  "5B9BD5"  # accent1
Set-ThemeColor $tcs 6  "ED7D31"  # accent2
Set-ThemeColor $tcs 7  "A5A5A5"  # accent3
Set-ThemeColor $tcs 8  "FFC000"  # accent4
Set-ThemeColor $tcs 9  "4472C4"  # accent5
Set-ThemeColor $tcs 10 "70AD47"  # accent6
Set-ThemeColor $tcs 11 "0563C1"  # hlink
Set-ThemeColor $tcs 12 "954F72"  # folHlink
